# EPBDS-1502 Added several tests.
#
# The original sheet contains one OpenL "Rules" test table (rows 4-10) and
# one matching "Testmethod" table (rows 13-17) describing a single test
# driverRiskScoreTest/driverRiskTest. This change renames that existing
# pair of tables to the "...1" variant and appends a second, similar pair
# of tables ("...2" variant, rows 20-25 and 28-31) that exercises an
# additional (empty) error case.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Rename the existing (first) block's titles to the "...1" variants.
# ---------------------------------------------------------------------
$ws.Range("B4").Value2  = "Rules DoubleValue driverRiskScoreTest1(String driverRisk)"
$ws.Range("B13").Value2 = "Testmethod driverRiskScoreTest1 driverRiskTest1"

# ---------------------------------------------------------------------
# 2. Duplicate the existing rows (content + formatting) into the new
#    block, cell by cell (keeps the original cell styles intact instead
#    of generating new style entries).
# ---------------------------------------------------------------------
$copyPairs = @(
    @("B4","B20"),  @("C4","C20"),  @("D4","D20"),
    @("B5","B21"),  @("C5","C21"),  @("D5","D21"),
    @("B6","B22"),  @("C6","C22"),  @("D6","D22"),
    @("B7","B23"),  @("C7","C23"),  @("D7","D23"),
    @("B8","B24"),  @("C8","C24"),  @("D8","D24"),
    @("B10","B25"), @("C10","C25"), @("D10","D25"),
    @("B13","B28"), @("C13","C28"), @("D13","D28"),
    @("B14","B29"), @("C14","C29"), @("D14","D29"),
    @("B15","B30"), @("C15","C30"), @("D15","D30"),
    @("B17","B31"), @("C17","C31"), @("D17","D31")
)

foreach ($pair in $copyPairs) {
    $ws.Range($pair[0]).Copy($ws.Range($pair[1]))
}

# ---------------------------------------------------------------------
# 3. Fix up the text that differs from a straight copy of block 1/2.
# ---------------------------------------------------------------------
$ws.Range("B20").Value2 = "Rules DoubleValue driverRiskScoreTest2(String driverRisk)"
$ws.Range("D25").Value2 = "'" + '=error("");0'
$ws.Range("B28").Value2 = "Testmethod driverRiskScoreTest2 driverRiskTest2"

# The new test expects an empty error (no exception text raised).
$ws.Range("D31").ClearContents()

# ---------------------------------------------------------------------
# 4. Re-create the merged cells for the two new blocks (mirroring the
#    merges already present in the first two blocks).
# ---------------------------------------------------------------------
$newMerges = @("B20:D20","B21:C21","B22:C22","B23:C23","B24:C24","B25:C25","B28:D28")
foreach ($m in $newMerges) {
    $ws.Range($m).Merge() | Out-Null
}

# ---------------------------------------------------------------------
# 5. Add the same dropdown validation used on B10 to the new B25 cell.
# ---------------------------------------------------------------------
$ws.Range("B25").Validation.Add(3, 1, 1, "driver_type") | Out-Null

# ---------------------------------------------------------------------
# 6. Update the selected cell to match the new active cell.
# ---------------------------------------------------------------------
$ws.Range("D31").Select() | Out-Null
